$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7
$ws.Cells.Item($row, 1).Value = 42611.884664351855
$ws.Cells.Item($row, 2).Value = 7
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Random"
